# LearningTree.xlsx - "Latest Changes and also added the Report Card Functionality"
#
# The STAGE and LMSProd sheets each keep a "current course" pointer in
# columns K/L of row 2 (CourseDesigner.../LearningCourse... test-data
# names). This commit rolls both sheets forward to a newer
# designer/course pair that was appended to the underlying test-data
# pool.

$wb = $excel.ActiveWorkbook

# ---- STAGE sheet (row 2, K:L) ---------------------------------------
$stage = $wb.Worksheets.Item("STAGE")

$stage.Range("K2").Value = "CourseDesigner86408"
$stage.Range("L2").Value = "LearningCourse26984"

# The refreshed pair also picks up the bordered/centered "active record"
# look (thin border on the right + bottom edges, centered both ways).
foreach ($addr in @("K2", "L2")) {
    $cell = $stage.Range($addr)
    $cell.Borders.Item(10).LineStyle = 1   # xlEdgeRight
    $cell.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
    $cell.HorizontalAlignment = -4108      # xlCenter
    $cell.VerticalAlignment = -4108        # xlCenter
}

# ---- LMSProd sheet (row 2, K:L) --------------------------------------
$lmsProd = $wb.Worksheets.Item("LMSProd")

$lmsProd.Range("K2").Value = "CourseDesigner23692"
$lmsProd.Range("L2").Value = "LearningCourse18179"
